$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Query1")

$ws.Range("C8").Value = "dsa"

$ws.Range("C8").Select()
